$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Populations")
$ws.Range("C6").Value = "M 0-14"
$ws.Range("C7").Value = "F 0-14"
$ws.Range("C8").Value = "M 15+"
$ws.Range("C9").Value = "F 15+"

$ws.Activate()
$null = $ws.Range("D13").Select()
